# "Better filtering for TEST and ERROR results"
#
# The workbook holds one row per sample with several pXRF element
# readings (columns B:K) normalized to [0,1]. One of the rows (A84,
# "TEST %") turned out to be a test/placeholder record rather than a
# real sample, so it is removed entirely (the rows below it shift up
# by one). Because the "Fe" column (G) is normalized against the full
# set of samples, removing that row changes the normalization range,
# so every remaining row's G value is rescaled with the same linear
# transform (derived from the reference data):
#   new_G = a * old_G + b

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$a = 1.0046342759312832
$b = -0.004634275931282797

# Find the row whose "Sample" column (A) holds the placeholder "TEST %"
# record and remove it entirely, shifting the rows below it up.
$lastRow = $ws.UsedRange.Rows.Count
$testRow = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2
    if ($label -eq "TEST %") {
        $testRow = $r
        break
    }
}

if ($testRow -gt 0) {
    $ws.Rows.Item($testRow).Delete()
}

# Recompute the "Fe" (column G) normalized values for every remaining
# data row using the updated normalization range.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $old = $cell.Value2
    if ($old -ne $null) {
        $cell.Value2 = $a * $old + $b
    }
}
